$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A65").Copy($ws.Range("A66"))
$ws.Range("A66").Value = 45970
$ws.Range("B66").Value = "15,2818"
$ws.Range("C66").Value = "15,8667"
$ws.Range("D66").Value = "15,2818"
$ws.Range("E66").Value = "15,2818"
